$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates per diff. Numeric-looking text values (e.g. "213.08") are
# written with a temporary text NumberFormat so Excel keeps them as strings
# (matching the source inlineStr cells) instead of auto-converting to numbers;
# the format is reset to General/Normal afterwards so no visible formatting changes.

$ws.Range("D2").Value = "26.166.50"
$ws.Range("E2").Value = "  +3.79%  "
$ws.Range("D3").Value = "1.602.60"
$ws.Range("E3").Value = "  +3.22%  "
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.51%  "
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.486"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.74%  "
$ws.Range("E8").Value = "  +3.82%  "
$ws.Range("E9").Value = "  +2.65%  "
$ws.Range("E10").Value = "  +2.21%  "
$ws.Range("E11").Value = "  +5.02%  "
$ws.Range("D12").Value = "1.825.94"
$ws.Range("E12").Value = "  +3.05%  "
$ws.Range("D13").Value = "1.603.92"
$ws.Range("E13").Value = "  +3.24%  "
$ws.Range("E14").Value = "  +1.35%  "
$ws.Range("E15").Value = "  +2.68%  "
$ws.Range("D16").Value = "26.151.99"
$ws.Range("E16").Value = "  +3.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.54"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.70%  "
$ws.Range("D18").Value = "0.0₃0720"
$ws.Range("E18").Value = "  +2.53%  "
$ws.Range("E19").Value = "  -0.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "205.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +12.14%  "
$ws.Range("E21").Value = "  +4.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.52%  "
$ws.Range("E23").Value = "  +3.31%  "
$ws.Range("E24").Value = "  +10.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.12%  "
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("E27").Value = "  -2.17%  "
$ws.Range("E28").Value = "  +3.74%  "
$ws.Range("E29").Value = "  +1.36%  "
$ws.Range("E30").Value = "  +1.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0472"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.66%  "
$ws.Range("E32").Value = "  +4.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.98"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.17%  "
$ws.Range("E34").Value = "  +2.97%  "
$ws.Range("E35").Value = "  +1.70%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0164"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.59%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "1.114.76"
$ws.Range("E37").Value = "  +3.27%  "
$ws.Range("E38").Value = "  -0.26%  "
$ws.Range("E39").Value = "  +3.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.780"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.494"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.781"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.35%  "
$ws.Range("D43").Value = "1.737.31"
$ws.Range("E43").Value = "  +3.00%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "93.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.90%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.39%  "
$ws.Range("E46").Value = "  +6.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "53.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0504"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("E49").Value = "  +0.89%  "
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.51%  "
